$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The rnaDate column (A) values for all data rows (2-29) were corrected
# from "09.04.20" to "09.19.20".
$range = $ws.Range("A2:A29")
$range.NumberFormat = "@"
$range.Value = "09.19.20"
$range.ClearFormats()

# Reflect the selection left on the sheet after making the edit.
$ws.Range("A3:A29").Select()
